{"js": "// Minor fixes on exercises\n// 1. Add \"space before\" (400 twips = 20pt) to the paragraph that reads\n//    \"\u0421\u044a\u0441\u0442\u043e\u044f\u043d\u0438\u0435 \u043d\u0430 \u043d\u0430\u0441\u0442\u043e\u044f\u0449\u0438\u044f \u0443\u0447\u0435\u0431\u0435\u043d \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b:\" (keeps its existing after=0).\n// 2. Remove the stray \"_GoBack\" bookmark left over near the status image.\n// 3. Add \"space before\" (400 twips = 20pt) and a hanging indent\n//    (left 357 twips / hanging 357 twips = 17.85pt) to the \"\u041b\u044e\u0431\u0438\u043c\u043e \u043c\u044f\u0441\u0442\u043e\"\n//    Heading 2 paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraphs by their text content so the script is resilient to\n// index shifts.\nlet statusParagraph = null;\nlet favoritePlaceHeading = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (statusParagraph === null && text.indexOf(\"\u0421\u044a\u0441\u0442\u043e\u044f\u043d\u0438\u0435 \u043d\u0430 \u043d\u0430\u0441\u0442\u043e\u044f\u0449\u0438\u044f \u0443\u0447\u0435\u0431\u0435\u043d \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\") !== -1) {\n    statusParagraph = paragraphs.items[i];\n  }\n  if (favoritePlaceHeading === null && text.indexOf(\"\u041b\u044e\u0431\u0438\u043c\u043e \u043c\u044f\u0441\u0442\u043e\") !== -1) {\n    favoritePlaceHeading = paragraphs.items[i];\n  }\n}\n\nif (statusParagraph) {\n  statusParagraph.paragraphFormat.spaceBefore = 20; // 400 twips\n}\n\nif (favoritePlaceHeading) {\n  favoritePlaceHeading.paragraphFormat.spaceBefore = 20; // 400 twips\n  favoritePlaceHeading.paragraphFormat.leftIndent = 17.85; // 357 twips\n  favoritePlaceHeading.paragraphFormat.firstLineIndent = -17.85; // hanging 357 twips\n}\n\n// Remove the leftover \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Minor fixes on exercises\n# 1. Add \"space before\" (400 twips = 20pt) to the paragraph that reads\n#    \"\u0421\u044a\u0441\u0442\u043e\u044f\u043d\u0438\u0435 \u043d\u0430 \u043d\u0430\u0441\u0442\u043e\u044f\u0449\u0438\u044f \u0443\u0447\u0435\u0431\u0435\u043d \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b:\" (keeps its existing after=0).\n# 2. Remove the stray \"_GoBack\" bookmark left over near the status image.\n# 3. Add \"space before\" (400 twips = 20pt) and a hanging indent\n#    (left 357 twips / hanging 357 twips = 17.85pt) to the \"\u041b\u044e\u0431\u0438\u043c\u043e \u043c\u044f\u0441\u0442\u043e\"\n#    Heading 2 paragraph.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n\n  if ($t -like \"*\u0421\u044a\u0441\u0442\u043e\u044f\u043d\u0438\u0435 \u043d\u0430 \u043d\u0430\u0441\u0442\u043e\u044f\u0449\u0438\u044f \u0443\u0447\u0435\u0431\u0435\u043d \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b*\") {\n    $p.Range.ParagraphFormat.SpaceBefore = 20\n  }\n\n  if ($t -like \"*\u041b\u044e\u0431\u0438\u043c\u043e \u043c\u044f\u0441\u0442\u043e*\") {\n    $p.Range.ParagraphFormat.SpaceBefore = 20\n    $p.Range.ParagraphFormat.LeftIndent = 17.85\n    $p.Range.ParagraphFormat.FirstLineIndent = -17.85\n  }\n}\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
